$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008 and 2009 rows (original rows 2 and 3), shifting the
# 2010 / 2011 rows up so they become rows 2 and 3.
$ws.Range("A2:K3").Delete()
